$d = $word.ActiveDocument

$replacements = @(
    @("246÷9=", "715÷7="),
    @("264÷7=", "645÷7="),
    @("417÷4=", "638÷6="),
    @("735÷8=", "203÷7="),
    @("501÷7=", "737÷2="),
    @("793÷5=", "419÷5="),
    @("813÷4=", "866÷4="),
    @("892÷7=", "384÷7="),
    @("303÷4=", "573÷5="),
    @("377÷8=", "653÷9="),
    @("494÷8=", "511÷9="),
    @("526÷7=", "635÷4="),
    @("673÷2=", "964÷8="),
    @("138÷4=", "355÷6="),
    @("554÷2=", "871÷3="),
    @("807÷3=", "651÷9="),
    @("797÷7=", "412÷7="),
    @("166÷5=", "664÷8="),
    @("857÷7=", "293÷4="),
    @("379÷5=", "670÷9="),
    @("333÷5=", "424÷5="),
    @("721÷4=", "323÷6="),
    @("370÷3=", "704÷5="),
    @("753÷7=", "608÷9="),
    @("108÷2=", "560÷3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
